$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D2 - User_Development row: append Communication_&_PR
$ws.Range("D2").Value = "UserDevelopment_Team_Denmark_Company,UserDevelopment_Team_Netherlands_Company,UserDevelopment_Team_Sweden_Company,UserDevelopment_Team_Germany_Company,Public_Relationship_Team,MarComm_Team,Event,Content_Team,Europe_UserDevelopment,Communication_&_PR"

# D3 - Sales_Operation row: append Europe_Business_Operation,Fleet_&_Business_Sales,Partner_Strategy,Retail_Sales,Sales_Planning
$ws.Range("D3").Value = "UserTeam_Norway_Company,UserOperations_Team_Germany_Company,UserOperations_Team_Netherlands_Company,UserOperations_Team_Denmark_Company,UserOperations_Team_Sweden_Company,Fleet_Planning_Team,Fleet_Operation_Team,Fleet_Management_Team_Netherlands_Company,Fleet_Management_Team_Sweden_Company,Fleet_Management_Team_Denmark_Company,Business_Development,Business_Intelligence,KA_Management_Team,Commercial_Product,Europe_Commercial_Operation,Europe_UserOperation_Department,Europe_Business_Operation,Fleet_&_Business_Sales,Partner_Strategy,Retail_Sales,Sales_Planning"

# D5 - Network_Infrastructure row: append EPX_PMO
$ws.Range("D5").Value = "Network_Development,Construction_Management,Cost_Management,Design_Management,NIO_House_Operation,PMO_Infrastructure,Europe_Space_Experience,EPX_PMO"

# D6 - Power_Operation row: append Power_Operation,Power Market_Launch & Enabling Team
$ws.Range("D6").Value = "Power_Business_Operation_Team,Power_Management_Team_Germany_Company,Power_Management_Team_Netherlands_Company,Strategy & Business_Development Team,Market_Launch & Enabling Team,Power_Management_Team_Sweden_Company,Power_Management_Team_Denmark_Company,Europe_Power_Operation_Department,Power_Operation,Power Market_Launch & Enabling Team"

# D7 - Service_Operation row: append Operation_Support,Parts_&_Logistics
$ws.Range("D7").Value = "Service_PMO,Service_Planning_Team,Service_Quality_Team,Spare_Parts_Team,Service_Operations_Team_Germany_Company,Service_Team_Norway_Company,Service_Operations_Team_Denmark_Company,Service_Operations_Team_Netherlands_Company,Service_Operations_Team_Sweden_Company,Europe_Service_Operation_Department,Operation_Support,Parts_&_Logistics"

# D12 - Backend_Support_Unit row: append several more departments
$ws.Range("D12").Value = "Europe_Business_HRBP_Department,Controlling_and_Planning_Department,Legal_EU_Department,NIO_Life,EU_Purchase,Operational_Procurement,Planning_Department,EU_Legal_Department,EHS,Digital_Development_PMO_Team,EU_PMK,(Digital) Sales Product Group,Product Marketing Department (PMK),Europe_Product_Experience_Department"
